$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Ningamma / 3GN22CS059, group 4 (Intelligent Learning Analytics Platform / Prof. Rajshekhar)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Ningamma"
$ws.Range("C6").Value = "3GN22CS059"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "Intelligent Learning Analytics Platform"
$ws.Range("F6").Value = "Prof. Rajshekhar "
$ws.Range("G6").Value = 47
$ws.Range("H6").Value = 42
$ws.Range("I6").Value = 46

# Row 7: Kanaka / 3GN22CS037, group 1 (Autism Detection / Prof.Bhimrao Patil)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Kanaka"
$ws.Range("C7").Value = "3GN22CS037"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Autism Detection"
$ws.Range("F7").Value = "Prof.Bhimrao Patil"
$ws.Range("G7").Value = 45
$ws.Range("H7").Value = 46
$ws.Range("I7").Value = 50

# Row 8: Akshata / 3GN22CS006, group 3 (Human Activity Recognisation / Prof. Johnwesley)
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Akshata"
$ws.Range("C8").Value = "3GN22CS006"
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "Human Activity Recognisation"
$ws.Range("F8").Value = "Prof. Johnwesley"
$ws.Range("G8").Value = 35
$ws.Range("H8").Value = 49
$ws.Range("I8").Value = 47

$ws.Range("E8:I8").Select()
